# Add a new ToDo row (row 3) to the 운영체제 sheet, mirroring the layout
# of the existing row 2 (과목/할 일/마감 기한/실제 마감일/완료 여부/중요도).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2 into row 3 first so formatting/typing of the copied cells
# (과목 column, the two corrupted "마감 기한"/"실제 마감일" placeholder
# cells) is preserved exactly.
$ws.Range("A2:F2").Copy($ws.Range("A3:F3"))

# Now overwrite the cells that actually change for the new row.
$ws.Range("B3").Value2 = "ㅁㄴㅇ"
$ws.Range("E3").Value2 = "진행"

# 중요도 (F3) should be the text "3" (same textual style as existing
# F2 = "2"), so force text formatting before assigning the numeric-looking
# string.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "3"
